$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-04 Sunday" "2024-08-05 Monday"
Replace-Text "40×98=" "19×41="
Replace-Text "35×70=" "15×56="
Replace-Text "18×30=" "18×20="
Replace-Text "58×18=" "53×32="
Replace-Text "29×30=" "43×76="
Replace-Text "32×24=" "21×89="
Replace-Text "77×90=" "73×49="
Replace-Text "65×23=" "92×32="
Replace-Text "32×29=" "89×77="
Replace-Text "87×53=" "68×14="
Replace-Text "78×14=" "13×68="
Replace-Text "46×28=" "42×78="
Replace-Text "70×14=" "25×71="
Replace-Text "58×92=" "38×73="
Replace-Text "86×14=" "71×61="
Replace-Text "82×13=" "74×68="
Replace-Text "70×83=" "39×75="
Replace-Text "54×70=" "63×40="
Replace-Text "77×60=" "17×36="
Replace-Text "85×76=" "97×56="
Replace-Text "14×82=" "39×42="
Replace-Text "39×89=" "82×63="
Replace-Text "53×42=" "74×78="
Replace-Text "92×47=" "49×18="
Replace-Text "15×82=" "75×58="
